$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A17").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("B17").Value = "89bdc2f6-0e22-47a8-b4f2-b7b5696fc495"
$ws.Range("C17").Value = 73.5
$ws.Range("D17").Value = "2025-08-07 23:02:17"
$ws.Range("E17").Value = "Paid"

$ws.Range("A18").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("B18").Value = "0947da20-6ab3-444d-97b4-2aa9c1662a75"
$ws.Range("C18").Value = 182.9
$ws.Range("D18").Value = "2025-08-07 23:14:24"
$ws.Range("E18").Value = "Paid"
